$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = 0
$ws.Cells.Item(1,2).Value = 1
$ws.Cells.Item(1,3).Value = 2

$ws.Cells.Item(2,1).Value = "period_name"
$ws.Cells.Item(2,2).Value = "CZK"
$ws.Cells.Item(2,3).Value = "česká koruna"

$ws.Cells.Item(3,1).Value = "2004q1"
$ws.Cells.Item(3,2).Value = 3.076923076923077
$ws.Cells.Item(3,3).Value = 10

$ws.Cells.Item(4,1).Value = "2004q2"
$ws.Cells.Item(4,2).Value = 0
$ws.Cells.Item(4,3).Value = 0

$ws.Cells.Item(5,1).Value = "2004q3"
$ws.Cells.Item(5,2).Value = 0
$ws.Cells.Item(5,3).Value = 0

$ws.Cells.Item(6,1).Value = "2004q4"
$ws.Cells.Item(6,2).Value = 0
$ws.Cells.Item(6,3).Value = 0

$ws.Cells.Item(7,1).Value = "2005q1"
$ws.Cells.Item(7,2).Value = 0
$ws.Cells.Item(7,3).Value = 0

$ws.Cells.Item(8,1).Value = "2005q2"
$ws.Cells.Item(8,2).Value = 0
$ws.Cells.Item(8,3).Value = 0

$ws.Cells.Item(9,1).Value = "2005q3"
$ws.Cells.Item(9,2).Value = 0
$ws.Cells.Item(9,3).Value = 0

$ws.Cells.Item(10,1).Value = "2005q4"
$ws.Cells.Item(10,2).Value = 0
$ws.Cells.Item(10,3).Value = 0

$ws.Cells.Item(11,1).Value = "2006q1"
$ws.Cells.Item(11,2).Value = 0
$ws.Cells.Item(11,3).Value = 0

$ws.Cells.Item(12,1).Value = "2006q2"
$ws.Cells.Item(12,2).Value = 0
$ws.Cells.Item(12,3).Value = 0

$ws.Cells.Item(13,1).Value = "2006q3"
$ws.Cells.Item(13,2).Value = 0
$ws.Cells.Item(13,3).Value = 0

$ws.Cells.Item(14,1).Value = "2006q4"
$ws.Cells.Item(14,2).Value = 2.846153846153846
$ws.Cells.Item(14,3).Value = 0

$ws.Cells.Item(15,1).Value = "2007q1"
$ws.Cells.Item(15,2).Value = 21.76923076923077
$ws.Cells.Item(15,3).Value = 0

$ws.Cells.Item(16,1).Value = "2007q2"
$ws.Cells.Item(16,2).Value = 21.30769230769231
$ws.Cells.Item(16,3).Value = 0

$ws.Cells.Item(17,1).Value = "2007q3"
$ws.Cells.Item(17,2).Value = 25.38461538461538
$ws.Cells.Item(17,3).Value = 0

$ws.Cells.Item(18,1).Value = "2007q4"
$ws.Cells.Item(18,2).Value = 25.53846153846154
$ws.Cells.Item(18,3).Value = 0

$ws.Cells.Item(19,1).Value = "2008q1"
$ws.Cells.Item(19,2).Value = 30.07692307692308
$ws.Cells.Item(19,3).Value = 0

$ws.Cells.Item(20,1).Value = "2008q2"
$ws.Cells.Item(20,2).Value = 37.53846153846154
$ws.Cells.Item(20,3).Value = 0

$ws.Cells.Item(21,1).Value = "2008q3"
$ws.Cells.Item(21,2).Value = 43.61538461538461
$ws.Cells.Item(21,3).Value = 0

$ws.Cells.Item(22,1).Value = "2008q4"
$ws.Cells.Item(22,2).Value = 37.15384615384615
$ws.Cells.Item(22,3).Value = 0

$ws.Cells.Item(23,1).Value = "2009q1"
$ws.Cells.Item(23,2).Value = 39.15384615384615
$ws.Cells.Item(23,3).Value = 0

$ws.Cells.Item(24,1).Value = "2009q2"
$ws.Cells.Item(24,2).Value = 37.46153846153846
$ws.Cells.Item(24,3).Value = 36

$ws.Cells.Item(25,1).Value = "2009q3"
$ws.Cells.Item(25,2).Value = 44.69230769230769
$ws.Cells.Item(25,3).Value = 46

$ws.Cells.Item(26,1).Value = "2009q4"
$ws.Cells.Item(26,2).Value = 44.15384615384615
$ws.Cells.Item(26,3).Value = 47.66666666666666

$ws.Cells.Item(27,1).Value = "2010q1"
$ws.Cells.Item(27,2).Value = 46.61538461538461
$ws.Cells.Item(27,3).Value = 60

$ws.Cells.Item(28,1).Value = "2010q2"
$ws.Cells.Item(28,2).Value = 47.69230769230769
$ws.Cells.Item(28,3).Value = 61.33333333333334

$ws.Cells.Item(29,1).Value = "2010q3"
$ws.Cells.Item(29,2).Value = 49.53846153846154
$ws.Cells.Item(29,3).Value = 39

$ws.Cells.Item(30,1).Value = "2010q4"
$ws.Cells.Item(30,2).Value = 45.61538461538461
$ws.Cells.Item(30,3).Value = 38.66666666666666

$ws.Cells.Item(31,1).Value = "2011q1"
$ws.Cells.Item(31,2).Value = 42.38461538461539
$ws.Cells.Item(31,3).Value = 29.33333333333333

$ws.Cells.Item(32,1).Value = "2011q2"
$ws.Cells.Item(32,2).Value = 46.23076923076923
$ws.Cells.Item(32,3).Value = 29

$ws.Cells.Item(33,1).Value = "2011q3"
$ws.Cells.Item(33,2).Value = 51.84615384615385
$ws.Cells.Item(33,3).Value = 27

$ws.Cells.Item(34,1).Value = "2011q4"
$ws.Cells.Item(34,2).Value = 47.61538461538461
$ws.Cells.Item(34,3).Value = 27.66666666666667

$ws.Cells.Item(35,1).Value = "2012q1"
$ws.Cells.Item(35,2).Value = 41.53846153846154
$ws.Cells.Item(35,3).Value = 23

$ws.Cells.Item(36,1).Value = "2012q2"
$ws.Cells.Item(36,2).Value = 45.15384615384615
$ws.Cells.Item(36,3).Value = 25.33333333333333

$ws.Cells.Item(37,1).Value = "2012q3"
$ws.Cells.Item(37,2).Value = 53
$ws.Cells.Item(37,3).Value = 29

$ws.Cells.Item(38,1).Value = "2012q4"
$ws.Cells.Item(38,2).Value = 45.15384615384615
$ws.Cells.Item(38,3).Value = 24

$ws.Cells.Item(39,1).Value = "2013q1"
$ws.Cells.Item(39,2).Value = 45
$ws.Cells.Item(39,3).Value = 21

$ws.Cells.Item(40,1).Value = "2013q2"
$ws.Cells.Item(40,2).Value = 47.38461538461539
$ws.Cells.Item(40,3).Value = 21.66666666666667

$ws.Cells.Item(41,1).Value = "2013q3"
$ws.Cells.Item(41,2).Value = 59.23076923076923
$ws.Cells.Item(41,3).Value = 31

$ws.Cells.Item(42,1).Value = "2013q4"
$ws.Cells.Item(42,2).Value = 57.84615384615385
$ws.Cells.Item(42,3).Value = 52

$ws.Cells.Item(43,1).Value = "2014q1"
$ws.Cells.Item(43,2).Value = 58.07692307692308
$ws.Cells.Item(43,3).Value = 24.33333333333333

$ws.Cells.Item(44,1).Value = "2014q2"
$ws.Cells.Item(44,2).Value = 59.53846153846154
$ws.Cells.Item(44,3).Value = 26.66666666666667

$ws.Cells.Item(45,1).Value = "2014q3"
$ws.Cells.Item(45,2).Value = 69.92307692307692
$ws.Cells.Item(45,3).Value = 26

$ws.Cells.Item(46,1).Value = "2014q4"
$ws.Cells.Item(46,2).Value = 68.15384615384616
$ws.Cells.Item(46,3).Value = 24

$ws.Cells.Item(47,1).Value = "2015q1"
$ws.Cells.Item(47,2).Value = 77.69230769230769
$ws.Cells.Item(47,3).Value = 25.66666666666667

$ws.Cells.Item(48,1).Value = "2015q2"
$ws.Cells.Item(48,2).Value = 80.07692307692308
$ws.Cells.Item(48,3).Value = 22.33333333333333

$ws.Cells.Item(49,1).Value = "2015q3"
$ws.Cells.Item(49,2).Value = 94.53846153846153
$ws.Cells.Item(49,3).Value = 24
